$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint 1")

$ws.Range("F9").Value = 0
$ws.Range("H9").ClearContents()

$ws.Range("C10").Value = 9
$ws.Range("D10").Value = 6
$ws.Range("E10").Value = 4
$ws.Range("F10").Value = 0

$ws.Range("C11").Value = 9
$ws.Range("D11").Value = 6
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 0
